$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data refreshed (new prices/scores) and the table was
# re-sorted descending by 최종점수 (final score), which swapped the
# 058470.KS row and the SamsungElec(005930.KS) row (rows 2 and 3).
# All other rows (4-7) keep their ticker but get refreshed metrics.

# Row 2: now 058470.KS
$ws.Range("B2").Value = "058470.KS,0P0000ASU1,98886"
$ws.Range("C2").Value = "058470.KS"
$ws.Range("D2").Value = 65100
$ws.Range("E2").Value = 64.7
$ws.Range("F2").Value = 0.46
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 63.7
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 54.85170003294819
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3: now SamsungElec / 005930.KS
$ws.Range("B3").Value = "SamsungElec"
$ws.Range("C3").Value = "005930.KS"
$ws.Range("D3").Value = 107400
$ws.Range("E3").Value = 60.4
$ws.Range("F3").Value = 6.87
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 58.5
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 54.85170003294819
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 4: 403870.KS (keep ticker, refreshed metrics)
$ws.Range("B4").Value = "403870.KS,0P0001PE9K,566428"
$ws.Range("C4").Value = "403870.KS"
$ws.Range("D4").Value = 30550
$ws.Range("E4").Value = 47.2
$ws.Range("F4").Value = 6.08
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 51.7
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 54.85170003294819
$ws.Range("O4").Value = "⚪ 중립 구간"

# Row 5: SK hynix
$ws.Range("B5").Value = "SK hynix"
$ws.Range("C5").Value = "000660.KS"
$ws.Range("D5").Value = 537000
$ws.Range("E5").Value = 32.3
$ws.Range("F5").Value = 1.32
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 66
$ws.Range("J5").Value = 66
$ws.Range("K5").Value = 48.9
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 54.85170003294819
$ws.Range("O5").Value = "⚪ 중립 구간"

# Row 6: DB HiTek
$ws.Range("B6").Value = "DB HiTek"
$ws.Range("C6").Value = "000990.KS"
$ws.Range("D6").Value = 65000
$ws.Range("E6").Value = 34.7
$ws.Range("F6").Value = 2.2
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 46
$ws.Range("J6").Value = 63
$ws.Range("K6").Value = 40.9
$ws.Range("L6").Value = "Pattern"
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 54.85170003294819
$ws.Range("O6").Value = "⚪ 중립 구간"

# Row 7: 240810.KS
$ws.Range("B7").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C7").Value = "240810.KS"
$ws.Range("D7").Value = 61000
$ws.Range("E7").Value = 36.7
$ws.Range("F7").Value = 1.16
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 46
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 40.9
$ws.Range("L7").Value = "Pattern"
$ws.Range("M7").Value = "⛔ 관망하십시오."
$ws.Range("N7").Value = 54.85170003294819
$ws.Range("O7").Value = "⚪ 중립 구간"
